$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new record row at row 44 (pushes the existing rows 44-116 down
# to 45-117, carrying all of their data/formatting with them).
$ws.Rows.Item(44).Insert()

# Populate the newly inserted row 44 with the new daily price record.
$ws.Range("A44").Value = 5
$ws.Range("B44").Value = "Macroferia Regional de Talca"
$ws.Range("C44").Value = "Maule"
$ws.Range("D44").Value = 44469
$ws.Range("E44").Value = 7
$ws.Range("F44").Value = 100112017
$ws.Range("G44").Value = "Apio"
$ws.Range("H44").Value = "Americana (o)"
$ws.Range("I44").Value = "Primera"
$ws.Range("J44").Value = 500
$ws.Range("K44").Value = 7500
$ws.Range("L44").Value = 7500
$ws.Range("M44").Value = 7500
$ws.Range("N44").Value = "`$/docena de matas"
$ws.Range("O44").Value = "Provincia del Elquí"
$ws.Range("P44").Value = 1250
$ws.Range("Q44").Value = 6
$ws.Range("R44").Value = "Hortaliza"
